$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 109.4120958805084
$ws.Range("C2").Value = 2.075103769611216
$ws.Range("D2").Value = 0.9311996459960937
$ws.Range("E2").Value = 0.03385167935236814
$ws.Range("J2").Value = 0.855072463768116
$ws.Range("K2").Value = 0.7536231884057971
$ws.Range("L2").Value = 0.7826086956521739
$ws.Range("M2").Value = 0.7971014492753623
$ws.Range("N2").Value = 0.8106796116504854
$ws.Range("O2").Value = 0.7998170817503869
$ws.Range("P2").Value = 0.03349188197737363
$ws.Range("Q2").Value = 13
$ws.Range("B3").Value = 218.711506652832
$ws.Range("C3").Value = 4.721050304847468
$ws.Range("D3").Value = 0.9865203857421875
$ws.Range("E3").Value = 0.04102906847897542
$ws.Range("J3").Value = 0.8454106280193237
$ws.Range("K3").Value = 0.748792270531401
$ws.Range("L3").Value = 0.7826086956521739
$ws.Range("M3").Value = 0.7971014492753623
$ws.Range("N3").Value = 0.8203883495145631
$ws.Range("O3").Value = 0.7988602785985648
$ws.Range("P3").Value = 0.03287675692407679
$ws.Range("Q3").Value = 16
$ws.Range("B4").Value = 445.9277801513672
$ws.Range("C4").Value = 3.047555800167626
$ws.Range("D4").Value = 1.008610057830811
$ws.Range("E4").Value = 0.03622302688370224
$ws.Range("J4").Value = 0.8454106280193237
$ws.Range("K4").Value = 0.7536231884057971
$ws.Range("L4").Value = 0.7826086956521739
$ws.Range("M4").Value = 0.8019323671497585
$ws.Range("N4").Value = 0.8252427184466019
$ws.Range("O4").Value = 0.801763519534731
$ws.Range("P4").Value = 0.03206484096109407
$ws.Range("Q4").Value = 10
$ws.Range("B5").Value = 116.4675260066986
$ws.Range("C5").Value = 1.388248479870066
$ws.Range("D5").Value = 1.078215885162354
$ws.Range("E5").Value = 0.1458360914246661
$ws.Range("J5").Value = 0.855072463768116
$ws.Range("K5").Value = 0.7536231884057971
$ws.Range("L5").Value = 0.7826086956521739
$ws.Range("M5").Value = 0.7971014492753623
$ws.Range("N5").Value = 0.8106796116504854
$ws.Range("O5").Value = 0.7998170817503869
$ws.Range("P5").Value = 0.03349188197737363
$ws.Range("Q5").Value = 13
$ws.Range("B6").Value = 225.5756698131561
$ws.Range("C6").Value = 1.227104272048248
$ws.Range("D6").Value = 1.012839698791504
$ws.Range("E6").Value = 0.02131754620670427
$ws.Range("J6").Value = 0.8454106280193237
$ws.Range("K6").Value = 0.748792270531401
$ws.Range("L6").Value = 0.7826086956521739
$ws.Range("M6").Value = 0.7971014492753623
$ws.Range("N6").Value = 0.8203883495145631
$ws.Range("O6").Value = 0.7988602785985648
$ws.Range("P6").Value = 0.03287675692407679
$ws.Range("Q6").Value = 16
$ws.Range("B7").Value = 446.1190069675446
$ws.Range("C7").Value = 1.412375711578195
$ws.Range("D7").Value = 1.072241497039795
$ws.Range("E7").Value = 0.09390612441359175
$ws.Range("J7").Value = 0.8454106280193237
$ws.Range("K7").Value = 0.7536231884057971
$ws.Range("L7").Value = 0.7826086956521739
$ws.Range("M7").Value = 0.8019323671497585
$ws.Range("N7").Value = 0.8252427184466019
$ws.Range("O7").Value = 0.801763519534731
$ws.Range("P7").Value = 0.03206484096109407
$ws.Range("Q7").Value = 10
$ws.Range("B8").Value = 114.7552185058594
$ws.Range("C8").Value = 0.6005699031590791
$ws.Range("D8").Value = 1.185709571838379
$ws.Range("E8").Value = 0.1955866029293977
$ws.Range("J8").Value = 0.855072463768116
$ws.Range("K8").Value = 0.7536231884057971
$ws.Range("L8").Value = 0.7826086956521739
$ws.Range("M8").Value = 0.7971014492753623
$ws.Range("N8").Value = 0.8106796116504854
$ws.Range("O8").Value = 0.7998170817503869
$ws.Range("P8").Value = 0.03349188197737363
$ws.Range("Q8").Value = 13
$ws.Range("B9").Value = 225.1608543395996
$ws.Range("C9").Value = 2.563349509251708
$ws.Range("D9").Value = 1.112521696090698
$ws.Range("E9").Value = 0.1745332125562456
$ws.Range("J9").Value = 0.8454106280193237
$ws.Range("K9").Value = 0.748792270531401
$ws.Range("L9").Value = 0.7826086956521739
$ws.Range("M9").Value = 0.7971014492753623
$ws.Range("N9").Value = 0.8203883495145631
$ws.Range("O9").Value = 0.7988602785985648
$ws.Range("P9").Value = 0.03287675692407679
$ws.Range("Q9").Value = 16
$ws.Range("B10").Value = 446.0811876296997
$ws.Range("C10").Value = 2.285915848026151
$ws.Range("D10").Value = 1.04566159248352
$ws.Range("E10").Value = 0.1022878067581385
$ws.Range("J10").Value = 0.8454106280193237
$ws.Range("K10").Value = 0.7536231884057971
$ws.Range("L10").Value = 0.7826086956521739
$ws.Range("M10").Value = 0.8019323671497585
$ws.Range("N10").Value = 0.8252427184466019
$ws.Range("O10").Value = 0.801763519534731
$ws.Range("P10").Value = 0.03206484096109407
$ws.Range("Q10").Value = 10
$ws.Range("B11").Value = 183.1310437202454
$ws.Range("C11").Value = 0.9778294508201462
$ws.Range("D11").Value = 1.027941989898682
$ws.Range("E11").Value = 0.0532828136564896
$ws.Range("J11").Value = 0.8599033816425121
$ws.Range("K11").Value = 0.7922705314009661
$ws.Range("L11").Value = 0.7874396135265701
$ws.Range("M11").Value = 0.8164251207729468
$ws.Range("N11").Value = 0.8203883495145631
$ws.Range("O11").Value = 0.8152853993715118
$ws.Range("P11").Value = 0.02578051723295969
$ws.Range("Q11").Value = 4
$ws.Range("B12").Value = 353.17311668396
$ws.Range("C12").Value = 2.326197653631079
$ws.Range("D12").Value = 1.083196449279785
$ws.Range("E12").Value = 0.1312995649419392
$ws.Range("J12").Value = 0.8599033816425121
$ws.Range("K12").Value = 0.7632850241545893
$ws.Range("L12").Value = 0.7874396135265701
$ws.Range("M12").Value = 0.821256038647343
$ws.Range("N12").Value = 0.8349514563106796
$ws.Range("O12").Value = 0.8133671028563387
$ws.Range("P12").Value = 0.03427592208899548
$ws.Range("Q12").Value = 7
$ws.Range("B13").Value = 684.1472810745239
$ws.Range("C13").Value = 4.323668383046725
$ws.Range("D13").Value = 0.963016128540039
$ws.Range("E13").Value = 0.01375185711530633
$ws.Range("J13").Value = 0.855072463768116
$ws.Range("K13").Value = 0.7536231884057971
$ws.Range("L13").Value = 0.7922705314009661
$ws.Range("M13").Value = 0.8309178743961353
$ws.Range("N13").Value = 0.8446601941747572
$ws.Range("O13").Value = 0.8153088504291544
$ws.Range("P13").Value = 0.037475786843105
$ws.Range("B14").Value = 186.8704335212707
$ws.Range("C14").Value = 4.647602261651786
$ws.Range("D14").Value = 0.9684895515441895
$ws.Range("E14").Value = 0.0159594481368628
$ws.Range("J14").Value = 0.8599033816425121
$ws.Range("K14").Value = 0.7922705314009661
$ws.Range("L14").Value = 0.7874396135265701
$ws.Range("M14").Value = 0.8164251207729468
$ws.Range("N14").Value = 0.8203883495145631
$ws.Range("O14").Value = 0.8152853993715118
$ws.Range("P14").Value = 0.02578051723295969
$ws.Range("Q14").Value = 4
$ws.Range("B15").Value = 355.2861651420593
$ws.Range("C15").Value = 1.286993147672724
$ws.Range("D15").Value = 0.9823958396911621
$ws.Range("E15").Value = 0.02165022199349082
$ws.Range("J15").Value = 0.8599033816425121
$ws.Range("K15").Value = 0.7632850241545893
$ws.Range("L15").Value = 0.7874396135265701
$ws.Range("M15").Value = 0.821256038647343
$ws.Range("N15").Value = 0.8349514563106796
$ws.Range("O15").Value = 0.8133671028563387
$ws.Range("P15").Value = 0.03427592208899548
$ws.Range("Q15").Value = 7
$ws.Range("B16").Value = 676.345999956131
$ws.Range("C16").Value = 11.49383885528823
$ws.Range("D16").Value = 1.049323320388794
$ws.Range("E16").Value = 0.1214148004774515
$ws.Range("J16").Value = 0.855072463768116
$ws.Range("K16").Value = 0.7536231884057971
$ws.Range("L16").Value = 0.7922705314009661
$ws.Range("M16").Value = 0.8309178743961353
$ws.Range("N16").Value = 0.8446601941747572
$ws.Range("O16").Value = 0.8153088504291544
$ws.Range("P16").Value = 0.037475786843105
$ws.Range("B17").Value = 183.172395324707
$ws.Range("C17").Value = 1.72310380731037
$ws.Range("D17").Value = 1.027873516082764
$ws.Range("E17").Value = 0.1031245211889095
$ws.Range("J17").Value = 0.8599033816425121
$ws.Range("K17").Value = 0.7922705314009661
$ws.Range("L17").Value = 0.7874396135265701
$ws.Range("M17").Value = 0.8164251207729468
$ws.Range("N17").Value = 0.8203883495145631
$ws.Range("O17").Value = 0.8152853993715118
$ws.Range("P17").Value = 0.02578051723295969
$ws.Range("Q17").Value = 4
$ws.Range("B18").Value = 361.0277002811432
$ws.Range("C18").Value = 6.55917342525869
$ws.Range("D18").Value = 0.9645584106445313
$ws.Range("E18").Value = 0.0491071801438441
$ws.Range("J18").Value = 0.8599033816425121
$ws.Range("K18").Value = 0.7632850241545893
$ws.Range("L18").Value = 0.7874396135265701
$ws.Range("M18").Value = 0.821256038647343
$ws.Range("N18").Value = 0.8349514563106796
$ws.Range("O18").Value = 0.8133671028563387
$ws.Range("P18").Value = 0.03427592208899548
$ws.Range("Q18").Value = 7
$ws.Range("B19").Value = 543.0735139369965
$ws.Range("C19").Value = 26.98052568034204
$ws.Range("D19").Value = 0.5907190799713135
$ws.Range("E19").Value = 0.09386576286604117
$ws.Range("J19").Value = 0.855072463768116
$ws.Range("K19").Value = 0.7536231884057971
$ws.Range("L19").Value = 0.7922705314009661
$ws.Range("M19").Value = 0.8309178743961353
$ws.Range("N19").Value = 0.8446601941747572
$ws.Range("O19").Value = 0.8153088504291544
$ws.Range("P19").Value = 0.037475786843105
